# Sprint 40 test-case report update
# - Fills in the "Day 9" summary numbers (Total testcase Written / Total
#   Execution / Total Review) which were previously left blank.
# - Leaves the scroll position / active selection on the sheet the way the
#   author left it after entering the numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enter the new "Day 9" totals (rows 51-53, column C) -------------------
$ws.Range("C51").Value = 2979
$ws.Range("C52").Value = 1455
$ws.Range("C53").Value = 852

# --- Nudge the merged header cells for the untouched day sections ----------
# Re-merging them (without changing anything visually) mirrors the way Excel
# re-emits the <mergeCells> list after such an edit, with the sections that
# were not touched directly being moved to the end of the list.
foreach ($ref in @("B2:C2", "B8:C8", "B14:C14", "B20:C20", "B26:C26")) {
    $ws.Range($ref).UnMerge()
    $ws.Range($ref).Merge()
}

# --- Update the window scroll position and current selection ---------------
$win = $excel.ActiveWindow
$win.ScrollRow = 42
$win.ScrollColumn = 1
$ws.Range("F49").Select()
